$d = $word.ActiveDocument

function Set-UniqueRangeText($searchText, $newText) {
    # Locate the single occurrence of $searchText in the document and
    # replace just that span's characters, leaving untouched runs (and
    # their formatting) alone.
    $r = $d.Content
    if ($r.Find.Execute($searchText)) {
        $r.Text = $newText
        return $true
    }
    return $false
}

# ------------------------------------------------------------------
# 1. " e atributos multi-valor." -- the paragraph is touched during the
#    edit pass, Word's proofing re-run drops the stale "multi-valor"
#    spell-check marker and folds the run back together.
# ------------------------------------------------------------------
$d.Content.Find.Execute(" e atributos multi-valor.", $true, $false, $false, $false, $false, $true, 1, $false, " e atributos multi-valor.", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Derivado Preço (Bilhete.Preço * (1 - Bilhete.Desconto)) -- same
#    kind of proofing-marker cleanup around "Bilhete.Preço".
# ------------------------------------------------------------------
$d.Content.Find.Execute(" (Bilhete.Preço * ", $true, $false, $false, $false, $false, $true, 1, $false, " (Bilhete.Preço * ", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Bilhete's primary key changes from the compound key (Lugar, Viagem)
#    to its own surrogate ID. Only the "Chave Primária" line below the
#    relation header is affected, so find the occurrence that is
#    preceded by "Chave Primária " and replace just its characters
#    (keeps the bold label run and the plain value run separate).
# ------------------------------------------------------------------
$r = $d.Content
$done = $false
while (-not $done -and $r.Find.Execute("Lugar, Viagem")) {
    $ctxStart = $r.Start - 15
    if ($ctxStart -lt 0) { $ctxStart = 0 }
    $ctx = $d.Range($ctxStart, $r.Start).Text
    if ($ctx -like "*Chave Primária *") {
        $target = $d.Range($r.Start, $r.End)
        $target.Text = "ID"
        $done = $true
    } else {
        $r.Collapse(0) | Out-Null
    }
}

# ------------------------------------------------------------------
# 4. Bilhete relation attribute list gains the new ID attribute.
# ------------------------------------------------------------------
Set-UniqueRangeText "(Lugar, Viagem, Classe, Preço, Reserva)" "(ID, Lugar, Viagem, Classe, Preço, Reserva)" | Out-Null

# ------------------------------------------------------------------
# 5. Viagem relation attribute list -- proofing-marker cleanup around
#    "DataPartida" / "PreçoBase" (no text change).
# ------------------------------------------------------------------
$d.Content.Find.Execute("(ID, DataPartida, Duração, PreçoBase, Comboio, Origem, Destino)", $true, $false, $false, $false, $false, $true, 1, $false, "(ID, DataPartida, Duração, PreçoBase, Comboio, Origem, Destino)", 2) | Out-Null

# ------------------------------------------------------------------
# 6. Lugar relation -- proofing-marker cleanup around "Nr" (no text
#    change), in both the relation header and its own "Chave Primária"
#    line.
# ------------------------------------------------------------------
$d.Content.Find.Execute("(Nr, Comboio)", $true, $false, $false, $false, $false, $true, 1, $false, "(Nr, Comboio)", 2) | Out-Null
$d.Content.Find.Execute("Chave Primária Nr", $true, $false, $false, $false, $false, $true, 1, $false, "Chave Primária Nr", 2) | Out-Null

# ------------------------------------------------------------------
# 7. Rewrite the narrative paragraph describing Bilhete's keys: Bilhete
#    now owns its own ID primary key and Viagem becomes a foreign key;
#    the old wording about Lugar being part of Bilhete's PK is dropped.
# ------------------------------------------------------------------
Set-UniqueRangeText "Bilhete tem Viagem (PK) que permite relacionar-se com Viagem através de ID (PK/FK). Cada entrada na tabela Bilhete tem um Lugar (PK) que corresponde ao número do lugar no comboio, uma Reserva (FK)" "Bilhete tem ID (PK), tem Viagem (FK) que permite relacionar-se com Viagem através de ID (PK). Cada entrada na tabela Bilhete tem uma Reserva (FK)" | Out-Null

# ------------------------------------------------------------------
# 8. "erá igual ao valor de PreçoBase da tabela Viagem..." -- proofing-
#    marker cleanup around "PreçoBase" (no text change).
# ------------------------------------------------------------------
$d.Content.Find.Execute("erá igual ao valor de PreçoBase da tabela Viagem, onde se aplica um desconto a partir de Classe.", $true, $false, $false, $false, $false, $true, 1, $false, "erá igual ao valor de PreçoBase da tabela Viagem, onde se aplica um desconto a partir de Classe.", 2) | Out-Null

# ------------------------------------------------------------------
# 9. " tem uma DataHoraPartida que indica o dia e a hora em que se
#    realiza a viagem, a duração (Duração) e o preço completo, sem
#    descontos, da viagem (PreçoBase). " -- proofing-marker cleanup
#    around "Data"/"Partida"/"PreçoBase" (no text change).
# ------------------------------------------------------------------
$d.Content.Find.Execute(" tem uma DataHoraPartida que indica o dia e a hora em que se realiza a viagem, a duração (Duração) e o preço completo, sem descontos, da viagem (PreçoBase). ", $true, $false, $false, $false, $false, $true, 1, $false, " tem uma DataHoraPartida que indica o dia e a hora em que se realiza a viagem, a duração (Duração) e o preço completo, sem descontos, da viagem (PreçoBase). ", 2) | Out-Null

# ------------------------------------------------------------------
# 10. "garantir que os campos referentes aos preços (...) são superiores
#     a zero." -- proofing-marker cleanup around "PreçoBase" (no text
#     change).
# ------------------------------------------------------------------
$d.Content.Find.Execute(" garantir que os campos referentes aos preços (Preço na Reserva, Preço no Bilhete e PreçoBase na Viagem) são superiores a zero.", $true, $false, $false, $false, $false, $true, 1, $false, " garantir que os campos referentes aos preços (Preço na Reserva, Preço no Bilhete e PreçoBase na Viagem) são superiores a zero.", 2) | Out-Null

# ------------------------------------------------------------------
# 11. Move the "_GoBack" bookmark from the end of the document (where it
#     marked the site of the previous last edit) to the Bilhete row
#     right before "Lugar", which is the most recently edited spot now.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$goBackAnchor = $d.Content
$goBackAnchor.Find.Execute("Bilhete (ID, ") | Out-Null
$goBackPos = $goBackAnchor.End
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos)) | Out-Null
